$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 183.22223
$ws.Range("I6").Value = 151.16667
$ws.Range("K6").Value = 453.50001
$ws.Range("M6").Value = -341.50001
$ws.Range("H9").Value = 290.44446
$ws.Range("I9").Value = 152.5
$ws.Range("J9").Value = 566.3333
$ws.Range("K9").Value = 152.5
$ws.Range("L9").Value = 566.3333
$ws.Range("M9").Value = 16.5
$ws.Range("N9").Value = -904.3333
$ws.Range("H43").Value = 5825.636
$ws.Range("J43").Value = 4022
$ws.Range("L43").Value = 4022
$ws.Range("N43").Value = -4160
$ws.Range("H62").Value = 6881.273
$ws.Range("I62").Value = 6724.875
$ws.Range("J62").Value = 7298.3335
$ws.Range("K62").Value = 6724.875
$ws.Range("L62").Value = 7298.3335
$ws.Range("M62").Value = -6100.875
$ws.Range("N62").Value = -8546.333500000001
$ws.Range("H65").Value = 6881.273
$ws.Range("I65").Value = 6724.875
$ws.Range("J65").Value = 7298.3335
$ws.Range("K65").Value = 33624.375
$ws.Range("L65").Value = 36491.6675
$ws.Range("M65").Value = -30504.375
$ws.Range("N65").Value = -42731.6675
$ws.Range("H103").Value = 41668144
$ws.Range("I103").Value = 799.5
$ws.Range("J103").Value = 50001612
$ws.Range("K103").Value = 2398.5
$ws.Range("L103").Value = 150004836
$ws.Range("M103").Value = -1812.5
$ws.Range("N103").Value = -150006008
$ws.Range("H112").Value = 4182.5835
$ws.Range("J112").Value = 5243.778
$ws.Range("L112").Value = 15731.334
$ws.Range("N112").Value = -17947.334
$ws.Range("H132").Value = 2530
$ws.Range("I132").Value = 2666.1738
$ws.Range("K132").Value = 7998.5214
$ws.Range("M132").Value = -5468.5214
$ws.Range("H135").Value = 2720.5715
$ws.Range("I135").Value = 261.25
$ws.Range("K135").Value = 2351.25
$ws.Range("M135").Value = 183.75
$ws.Range("H138").Value = 3188.275
$ws.Range("I138").Value = 1810.9231
$ws.Range("J138").Value = 3851.4443
$ws.Range("K138").Value = 5432.7693
$ws.Range("L138").Value = 11554.3329
$ws.Range("M138").Value = -292.7692999999999
$ws.Range("N138").Value = -21834.3329

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7002.7114
$ws.Range("I32").Value = 6335.551
$ws.Range("K32").Value = 6335.551
$ws.Range("M32").Value = -6048.551
$ws.Range("H45").Value = 9999.5
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 9999.5
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 9999.5
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -10753.5
$ws.Range("H61").Value = 9550571
$ws.Range("I61").Value = 10531856
$ws.Range("K61").Value = 10531856
$ws.Range("M61").Value = -10531644
$ws.Range("H102").Value = 8493.546
$ws.Range("J102").Value = 9684.200000000001
$ws.Range("L102").Value = 9684.200000000001
$ws.Range("N102").Value = -12928.2
$ws.Range("H110").Value = 3974.2856
$ws.Range("J110").Value = 4973.125
$ws.Range("L110").Value = 4973.125
$ws.Range("N110").Value = -9063.125
$ws.Range("H132").Value = 3450658.8
$ws.Range("I132").Value = 2316.75
$ws.Range("K132").Value = 6950.25
$ws.Range("M132").Value = -4420.25
$ws.Range("H136").Value = 9550571
$ws.Range("I136").Value = 10531856
$ws.Range("K136").Value = 31595568
$ws.Range("M136").Value = -31593018
$ws.Range("H140").Value = 249999.67
$ws.Range("J140").Value = 249999.67
$ws.Range("L140").Value = 249999.67
$ws.Range("N140").Value = -260359.67

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 656666.7
$ws.Range("I105").Value = 995819.6
$ws.Range("K105").Value = 995819.6
$ws.Range("M105").Value = -994072.6
$ws.Range("H107").Value = 3024.125
$ws.Range("I107").Value = 3205.2856
$ws.Range("J107").Value = 1756
$ws.Range("K107").Value = 3205.2856
$ws.Range("L107").Value = 1756
$ws.Range("M107").Value = -1285.2856
$ws.Range("N107").Value = -5596
$ws.Range("H109").Value = 99999
$ws.Range("J109").Value = 99999
$ws.Range("L109").Value = 99999
$ws.Range("N109").Value = -102773

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 40002704
$ws.Range("J31").Value = 3550
$ws.Range("L31").Value = 3550
$ws.Range("N31").Value = -4140
$ws.Range("H34").Value = 40002704
$ws.Range("J34").Value = 3550
$ws.Range("L34").Value = 3550
$ws.Range("N34").Value = -3954
$ws.Range("H107").Value = 1440.1111
$ws.Range("I107").Value = 326.41666
$ws.Range("J107").Value = 3667.5
$ws.Range("K107").Value = 326.41666
$ws.Range("L107").Value = 3667.5
$ws.Range("M107").Value = 1593.58334
$ws.Range("N107").Value = -7507.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 24436.916
$ws.Range("I63").Value = 6970.6665
$ws.Range("J63").Value = 30259
$ws.Range("K63").Value = 20911.9995
$ws.Range("L63").Value = 90777
$ws.Range("M63").Value = -20162.9995
$ws.Range("N63").Value = -92275
$ws.Range("H66").Value = 24436.916
$ws.Range("I66").Value = 6970.6665
$ws.Range("J66").Value = 30259
$ws.Range("K66").Value = 62735.9985
$ws.Range("L66").Value = 272331
$ws.Range("M66").Value = -58991.9985
$ws.Range("N66").Value = -279819
$ws.Range("H92").Value = 58
$ws.Range("J92").Value = 58
$ws.Range("L92").Value = 174
$ws.Range("N92").Value = -2670
$ws.Range("H114").Value = 11954.333
$ws.Range("J114").Value = 17916.5
$ws.Range("L114").Value = 53749.5
$ws.Range("N114").Value = -60257.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 1000
$ws.Range("K4").Value = 1000
$ws.Range("M4").Value = -888
$ws.Range("H5").Value = 19832.334
$ws.Range("I5").Value = 17248.5
$ws.Range("K5").Value = 17248.5
$ws.Range("M5").Value = -17136.5
$ws.Range("H80").Value = 2670.9092
$ws.Range("I80").Value = 2472.625
$ws.Range("K80").Value = 2472.625
$ws.Range("M80").Value = -1474.625
$ws.Range("H83").Value = 2670.9092
$ws.Range("I83").Value = 2472.625
$ws.Range("K83").Value = 12363.125
$ws.Range("M83").Value = -7371.125
$ws.Range("H102").Value = 2387.6
$ws.Range("I102").Value = 2486.2222
$ws.Range("K102").Value = 2486.2222
$ws.Range("M102").Value = -864.2222000000002
$ws.Range("H107").Value = 793.3333
$ws.Range("I107").Value = 190.5
$ws.Range("J107").Value = 1999
$ws.Range("K107").Value = 190.5
$ws.Range("L107").Value = 1999
$ws.Range("M107").Value = 1729.5
$ws.Range("N107").Value = -5839
$ws.Range("H123").Value = 98666.664
$ws.Range("J123").Value = 98666.664
$ws.Range("L123").Value = 98666.664
$ws.Range("N123").Value = -103566.664
$ws.Range("H132").Value = 6368008.5
$ws.Range("I132").Value = 4666.44
$ws.Range("K132").Value = 13999.32
$ws.Range("M132").Value = -11469.32

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6690.077
$ws.Range("I22").Value = 13332.667
$ws.Range("K22").Value = 13332.667
$ws.Range("M22").Value = -13037.667
$ws.Range("H27").Value = 6690.077
$ws.Range("I27").Value = 13332.667
$ws.Range("K27").Value = 13332.667
$ws.Range("M27").Value = -13225.667
$ws.Range("H40").Value = 6011.9443
$ws.Range("I40").Value = 6071.4707
$ws.Range("K40").Value = 6071.4707
$ws.Range("M40").Value = -5935.4707
$ws.Range("H93").Value = 4278312
$ws.Range("I93").Value = 5374.5
$ws.Range("J93").Value = 6177395.5
$ws.Range("K93").Value = 5374.5
$ws.Range("L93").Value = 6177395.5
$ws.Range("M93").Value = -4126.5
$ws.Range("N93").Value = -6179891.5
$ws.Range("H132").Value = 3137.0334
$ws.Range("I132").Value = 2076.9167
$ws.Range("J132").Value = 4727.2085
$ws.Range("K132").Value = 6230.750100000001
$ws.Range("L132").Value = 14181.6255
$ws.Range("M132").Value = -3700.750100000001
$ws.Range("N132").Value = -19241.6255

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1988.6154
$ws.Range("I81").Value = 2411.625
$ws.Range("J81").Value = 1311.8
$ws.Range("K81").Value = 4823.25
$ws.Range("L81").Value = 2623.6
$ws.Range("M81").Value = -3762.25
$ws.Range("N81").Value = -4745.6
$ws.Range("H84").Value = 1988.6154
$ws.Range("I84").Value = 2411.625
$ws.Range("J84").Value = 1311.8
$ws.Range("K84").Value = 24116.25
$ws.Range("L84").Value = 13118
$ws.Range("M84").Value = -18812.25
$ws.Range("N84").Value = -23726
$ws.Range("H107").Value = 5311.0625
$ws.Range("I107").Value = 3767.6924
$ws.Range("K107").Value = 11303.0772
$ws.Range("M107").Value = -9383.0772
$ws.Range("H122").Value = 2561.077
$ws.Range("I122").Value = 2328.1428
$ws.Range("J122").Value = 2832.8333
$ws.Range("K122").Value = 6984.428400000001
$ws.Range("L122").Value = 8498.499899999999
$ws.Range("M122").Value = -4534.428400000001
$ws.Range("N122").Value = -13398.4999
$ws.Range("H132").Value = 373617.47
$ws.Range("I132").Value = 3369.182
$ws.Range("K132").Value = 10107.546
$ws.Range("M132").Value = -7577.545999999998
